$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing date-formatted cell (F103) as a style template so the new
# date cells reuse the workbook's existing number-format style instead of
# Excel allocating a brand new custom numFmt.
$dateTemplate = $ws.Range("F103")

# --- Row 104 (ID 103): Clone Linked List with Random Pointer ---
$ws.Range("B104").Value = "Linked List"
$ws.Range("C104").Value = "Clone Linked List with Random Pointer"
$ws.Range("D104").Value = "Medium"
$ws.Range("E104").Value = "Done"
$dateTemplate.Copy($ws.Range("F104"))
$ws.Range("F104").Value = 45926
$ws.Range("G104").Value = "O(n)"
$ws.Range("H104").Value = "O(n)"
$ws.Range("I104").Value = "Using HashMap"

# --- Row 105 (ID 104): Palindrome Linked List ---
$ws.Range("B105").Value = "Linked List"
$ws.Range("C105").Value = "Palindrom Linked List"
$ws.Range("D105").Value = "Easy"
$ws.Range("E105").Value = "Done"
$dateTemplate.Copy($ws.Range("F105"))
$ws.Range("F105").Value = 45926
$ws.Range("G105").Value = "O(n)"
$ws.Range("H105").Value = "O(1)"
$ws.Range("I105").Value = "Reversing"

# --- Row 106 (ID 105): Merge Sort ---
$ws.Range("B106").Value = "Linked List"
$ws.Range("C106").Value = "Merge Sort"
$ws.Range("D106").Value = "Medium"
$ws.Range("E106").Value = "Done"
$dateTemplate.Copy($ws.Range("F106"))
$ws.Range("F106").Value = 45926
$ws.Range("G106").Value = "O(n * log n)"
$ws.Range("H106").Value = "O(log n)"

# Move the current selection from A107 to B107, matching the saved view state.
[void]$ws.Range("B107").Select()
